$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from the existing header cell (H1) onto the two new
# header cells so they get the same bold/centered/bordered formatting.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New column data (I0, IF)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 4

$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 5

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 4

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 3

$ws.Range("I7").Value = 3
$ws.Range("J7").Value = 4

$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 2
